$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts old C..I to D..J)
$ws.Columns("C:C").Insert()

# Header
$ws.Range("C1").Value = "Industry"

# Fill in Industry values for each holding row
$ws.Range("C2").Value = "Pharmaceuticals & Biotechnology"
$ws.Range("C3").Value = "Healthcare Services"
$ws.Range("C4").Value = "Consumer Durables"
$ws.Range("C5").Value = "Retailing"
$ws.Range("C6").Value = "Transport Services"
$ws.Range("C7").Value = "Healthcare Services"
$ws.Range("C8").Value = "Banks"
$ws.Range("C9").Value = "Capital Markets"
$ws.Range("C10").Value = "Cement & Cement Products"
$ws.Range("C11").Value = "Banks"
$ws.Range("C12").Value = "Banks"
$ws.Range("C13").Value = "Industrial Manufacturing"
$ws.Range("C14").Value = "Finance"
$ws.Range("C15").Value = "Retailing"
$ws.Range("C16").Value = "Pharmaceuticals & Biotechnology"
$ws.Range("C17").Value = "Transport Services"
$ws.Range("C18").Value = "Healthcare Services"
$ws.Range("C19").Value = "Finance"
$ws.Range("C20").Value = "Banks"
$ws.Range("C21").Value = "Commercial Services & Supplies"
$ws.Range("C22").Value = "Automobiles"
$ws.Range("C23").Value = "Auto Components"
$ws.Range("C24").Value = "Auto Components"
$ws.Range("C25").Value = "IT - Software"
$ws.Range("C26").Value = "Retailing"
$ws.Range("C27").Value = "Insurance"
$ws.Range("C28").Value = "Realty"
$ws.Range("C29").Value = "Healthcare Services"
$ws.Range("C30").Value = "Realty"
$ws.Range("C31").Value = "Paper, Forest & Jute Products"
$ws.Range("C32").Value = "Leisure Services"
$ws.Range("C33").Value = "Capital Markets"
$ws.Range("C34").Value = "Agricultural, Commercial & Construction Vehicles"
$ws.Range("C35").Value = "IT - Software"
$ws.Range("C36").Value = "Capital Markets"
$ws.Range("C37").Value = "Consumer Durables"
$ws.Range("C38").Value = "Industrial Manufacturing"
$ws.Range("C39").Value = "Consumer Durables"
$ws.Range("C40").Value = "Electrical Equipment"
$ws.Range("C41").Value = "Finance"
$ws.Range("C42").Value = "Industrial Products"
$ws.Range("C43").Value = "Electrical Equipment"
$ws.Range("C44").Value = "Food Products"
$ws.Range("C45").Value = "IT - Software"
$ws.Range("C46").Value = "Capital Markets"
$ws.Range("C47").Value = "Pharmaceuticals & Biotechnology"
$ws.Range("C48").Value = "Finance"
$ws.Range("C49").Value = "Finance"
$ws.Range("C50").Value = "Realty"
$ws.Range("C51").Value = "Industrial Products"
$ws.Range("C52").Value = "Construction"
$ws.Range("C53").Value = "Industrial Products"
$ws.Range("C54").Value = "Pharmaceuticals & Biotechnology"
$ws.Range("C55").Value = "Pharmaceuticals & Biotechnology"
$ws.Range("C56").Value = "Industrial Products"
$ws.Range("C57").Value = "Chemicals & Petrochemicals"
$ws.Range("C58").Value = "Agricultural Food & other Products"
$ws.Range("C59").Value = "Cement & Cement Products"
$ws.Range("C60").Value = "Chemicals & Petrochemicals"
$ws.Range("C61").Value = "Realty"
$ws.Range("C62").Value = "Industrial Products"
$ws.Range("C63").Value = "Healthcare Services"
$ws.Range("C64").Value = "Industrial Products"
$ws.Range("C65").Value = "Pharmaceuticals & Biotechnology"
$ws.Range("C66").Value = "Capital Markets"
$ws.Range("C67").Value = "Agricultural, Commercial & Construction Vehicles"
$ws.Range("C68").Value = "Capital Markets"
$ws.Range("C69").Value = "Capital Markets"
$ws.Range("C70").Value = "Healthcare Services"
$ws.Range("C71").Value = "Retailing"
$ws.Range("C72").Value = "Retailing"
$ws.Range("C73").Value = "Food Products"
